$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3021.5
$ws.Range("J17").Value = 5000
$ws.Range("L17").Value = 15000
$ws.Range("N17").Value = -15336

$ws.Range("H62").Value = 5853.5
$ws.Range("J62").Value = 5328.6665
$ws.Range("L62").Value = 5328.6665
$ws.Range("N62").Value = -6576.6665

$ws.Range("H64").Value = 14317.682
$ws.Range("I64").Value = 9582.916999999999
$ws.Range("J64").Value = 19999.4
$ws.Range("K64").Value = 9582.916999999999
$ws.Range("L64").Value = 19999.4
$ws.Range("M64").Value = -9334.916999999999
$ws.Range("N64").Value = -20495.4

$ws.Range("H65").Value = 5853.5
$ws.Range("J65").Value = 5328.6665
$ws.Range("L65").Value = 26643.3325
$ws.Range("N65").Value = -32883.3325

$ws.Range("H67").Value = 14317.682
$ws.Range("I67").Value = 9582.916999999999
$ws.Range("J67").Value = 19999.4
$ws.Range("K67").Value = 9582.916999999999
$ws.Range("L67").Value = 19999.4
$ws.Range("M67").Value = -8724.916999999999
$ws.Range("N67").Value = -21715.4

$ws.Range("H69").Value = 5979.8335
$ws.Range("I69").Value = 5949.5
$ws.Range("K69").Value = 17848.5
$ws.Range("M69").Value = -16974.5

$ws.Range("H72").Value = 5979.8335
$ws.Range("I72").Value = 5949.5
$ws.Range("K72").Value = 53545.5
$ws.Range("M72").Value = -49177.5

$ws.Range("H74").Value = 2990
$ws.Range("I74").Value = 2990
$ws.Range("K74").Value = 2990
$ws.Range("M74").Value = -2054

$ws.Range("H77").Value = 2990
$ws.Range("I77").Value = 2990
$ws.Range("K77").Value = 14950
$ws.Range("M77").Value = -10270

$ws.Range("H98").Value = 1046.25
$ws.Range("I98").Value = 1127.5
$ws.Range("K98").Value = 1127.5
$ws.Range("M98").Value = 370.5

$ws.Range("H116").Value = 4659.375
$ws.Range("I116").Value = 4255.2
$ws.Range("J116").Value = 5333
$ws.Range("K116").Value = 4255.2
$ws.Range("L116").Value = 5333
$ws.Range("M116").Value = -813.1999999999998
$ws.Range("N116").Value = -12217

$ws.Range("H122").Value = 1046.25
$ws.Range("I122").Value = 1127.5
$ws.Range("K122").Value = 3382.5
$ws.Range("M122").Value = -932.5

$ws.Range("H132").Value = 4439.6
$ws.Range("I132").Value = 3674.5
$ws.Range("K132").Value = 11023.5
$ws.Range("M132").Value = -8493.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H122").Value = 1516.6
$ws.Range("I122").Value = 1516.6
$ws.Range("K122").Value = 4549.799999999999
$ws.Range("M122").Value = -2099.799999999999

$ws.Range("H132").Value = 2485.8462
$ws.Range("I132").Value = 1887.875
$ws.Range("K132").Value = 5663.625
$ws.Range("M132").Value = -3133.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1881.6
$ws.Range("I94").Value = 545.2857
$ws.Range("K94").Value = 545.2857
$ws.Range("M94").Value = -94.28570000000002

$ws.Range("H105").Value = 2184.7273
$ws.Range("I105").Value = 2255.5
$ws.Range("K105").Value = 2255.5
$ws.Range("M105").Value = -508.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3185.7646
$ws.Range("I31").Value = 2939.9285
$ws.Range("K31").Value = 2939.9285
$ws.Range("M31").Value = -2644.9285

$ws.Range("H32").Value = 2832.8333
$ws.Range("I32").Value = 2899.4
$ws.Range("K32").Value = 2899.4
$ws.Range("M32").Value = -2583.4

$ws.Range("H34").Value = 3185.7646
$ws.Range("I34").Value = 2939.9285
$ws.Range("K34").Value = 2939.9285
$ws.Range("M34").Value = -2737.9285

$ws.Range("H58").Value = 9838.125
$ws.Range("I58").Value = 9813
$ws.Range("J58").Value = 10014
$ws.Range("K58").Value = 9813
$ws.Range("L58").Value = 10014
$ws.Range("M58").Value = -9610
$ws.Range("N58").Value = -10420

$ws.Range("H86").Value = 34854484
$ws.Range("I86").Value = 34854484
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 34854484
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -34853361

$ws.Range("H89").Value = 34854484
$ws.Range("I89").Value = 34854484
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 174272420
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -174266804

$ws.Range("H99").Value = 6821.4707
$ws.Range("I99").Value = 6186.5835
$ws.Range("K99").Value = 6186.5835
$ws.Range("M99").Value = -4688.5835

$ws.Range("H103").Value = 8000
$ws.Range("I103").Value = 8000
$ws.Range("K103").Value = 8000
$ws.Range("M103").Value = -6828

$ws.Range("H107").Value = 367.27274
$ws.Range("I107").Value = 205.66667
$ws.Range("K107").Value = 205.66667
$ws.Range("M107").Value = 1714.33333

$ws.Range("H126").Value = 6821.4707
$ws.Range("I126").Value = 6186.5835
$ws.Range("K126").Value = 18559.7505
$ws.Range("M126").Value = -16089.7505

$ws.Range("H132").Value = 2763.7368
$ws.Range("I132").Value = 1721.2858
$ws.Range("J132").Value = 5682.6
$ws.Range("K132").Value = 5163.857400000001
$ws.Range("L132").Value = 17047.8
$ws.Range("M132").Value = -2633.857400000001
$ws.Range("N132").Value = -22107.8

$ws.Range("H136").Value = 9838.125
$ws.Range("I136").Value = 9813
$ws.Range("J136").Value = 10014
$ws.Range("K136").Value = 29439
$ws.Range("L136").Value = 30042
$ws.Range("M136").Value = -26889
$ws.Range("N136").Value = -35142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 1699.8
$ws.Range("I87").Value = 1749.75
$ws.Range("J87").Value = 1500
$ws.Range("K87").Value = 5249.25
$ws.Range("L87").Value = 4500
$ws.Range("M87").Value = -4001.25
$ws.Range("N87").Value = -6996

$ws.Range("H90").Value = 1699.8
$ws.Range("I90").Value = 1749.75
$ws.Range("J90").Value = 1500
$ws.Range("K90").Value = 15747.75
$ws.Range("L90").Value = 13500
$ws.Range("M90").Value = -9507.75
$ws.Range("N90").Value = -25980

$ws.Range("H109").Value = 596.75
$ws.Range("I109").Value = 431.42856
$ws.Range("J109").Value = 828.2
$ws.Range("K109").Value = 1294.28568
$ws.Range("L109").Value = 2484.6
$ws.Range("M109").Value = -254.28568
$ws.Range("N109").Value = -4564.6

$ws.Range("H137").Value = 1615.2
$ws.Range("J137").Value = 3010.6667
$ws.Range("L137").Value = 9032.000100000001
$ws.Range("N137").Value = -19232.0001

$ws.Range("H139").Value = 2352.125
$ws.Range("I139").Value = 764.8
$ws.Range("K139").Value = 2294.4
$ws.Range("M139").Value = 2845.6

$ws.Range("H140").Value = 1000
$ws.Range("I140").Value = 1000
$ws.Range("K140").Value = 3000
$ws.Range("M140").Value = 2180

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 10418.333
$ws.Range("I36").Value = 10418.333
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 10418.333
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -9933.333000000001

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("N80").Value = 0

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("N83").Value = 0

$ws.Range("H132").Value = 3242.2307
$ws.Range("I132").Value = 2725.25
$ws.Range("J132").Value = 4069.4
$ws.Range("K132").Value = 8175.75
$ws.Range("L132").Value = 12208.2
$ws.Range("M132").Value = -5645.75
$ws.Range("N132").Value = -17268.2

$ws.Range("H136").Value = 33664.43
$ws.Range("J136").Value = 33664.43
$ws.Range("L136").Value = 100993.29
$ws.Range("N136").Value = -106093.29

$ws.Range("H141").Value = 45500
$ws.Range("J141").Value = 45500
$ws.Range("L141").Value = 45500
$ws.Range("N141").Value = -55860

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2903
$ws.Range("I46").Value = 899
$ws.Range("K46").Value = 899
$ws.Range("M46").Value = -711

$ws.Range("H93").Value = 874.75
$ws.Range("J93").Value = 800
$ws.Range("L93").Value = 800
$ws.Range("N93").Value = -3296

$ws.Range("H122").Value = 2681.6667
$ws.Range("I122").Value = 2681.6667
$ws.Range("K122").Value = 8045.000100000001
$ws.Range("M122").Value = -5595.000100000001

$ws.Range("H132").Value = 4457.773
$ws.Range("I132").Value = 3435.6875
$ws.Range("J132").Value = 7183.3335
$ws.Range("K132").Value = 10307.0625
$ws.Range("L132").Value = 21550.0005
$ws.Range("M132").Value = -7777.0625
$ws.Range("N132").Value = -26610.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4004
$ws.Range("I122").Value = 4004
$ws.Range("K122").Value = 12012
$ws.Range("M122").Value = -9562
